$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers to Excel's parser (single
# decimal point, no thousands separators) must be forced to Text format first,
# otherwise Excel auto-converts the assigned string into a numeric value and the
# literal text (e.g. "0.0000126") would be stored/rounded as a number instead.
$numericLooking = @("D5","D6","D7","D9","D11","D14","D15","D18","D20","D21","D22","D23","D24","D26","D27","D28","D34","D35","D37","D40","D41","D45","D46","D48")
foreach ($ref in $numericLooking) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated price values for the numeric-looking cells
$ws.Range("D5").Value = "400.75"
$ws.Range("D6").Value = "126.06"
$ws.Range("D7").Value = "0.589"
$ws.Range("D9").Value = "0.660"
$ws.Range("D11").Value = "40.88"
$ws.Range("D14").Value = "8.27"
$ws.Range("D15").Value = "19.29"
$ws.Range("D18").Value = "11.24"
$ws.Range("D20").Value = "0.0000126"
$ws.Range("D21").Value = "3.20"
$ws.Range("D22").Value = "79.88"
$ws.Range("D23").Value = "12.69"
$ws.Range("D24").Value = "298.60"
$ws.Range("D26").Value = "4.74"
$ws.Range("D27").Value = "29.01"
$ws.Range("D28").Value = "8.19"
$ws.Range("D34").Value = "2.49"
$ws.Range("D35").Value = "40.95"
$ws.Range("D37").Value = "52.02"
$ws.Range("D40").Value = "2.92"
$ws.Range("D41").Value = "137.10"
$ws.Range("D45").Value = "3.90"
$ws.Range("D46").Value = "16.55"
$ws.Range("D48").Value = "21.18"

# Restore the default (General) style on those cells so only the value changed,
# matching the original workbook formatting.
foreach ($ref in $numericLooking) {
    $ws.Range($ref).Style = "Normal"
}

# Remaining cell updates (multi-dot prices, URLs, names, percentage strings) are
# already unambiguous text to Excel's auto-detection and can be set directly.
$ws.Range("D2").Value = "61.551.87"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "3.354.66"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E6").Value = "  +7.65%  "
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +4.46%  "
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "3.881.44"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "3.346.04"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "61.423.29"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("E18").Value = "  +3.24%  "
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  +7.05%  "
$ws.Range("E21").Value = "  -4.81%  "
$ws.Range("E22").Value = "  +6.68%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  +11.18%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E28").Value = "  +7.63%  "
$ws.Range("E29").Value = "  -7.41%  "
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("E40").Value = "  -6.75%  "
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").Value = "3.684.56"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "2.095.84"
$ws.Range("E50").Value = "  -3.49%  "
$ws.Range("E51").Value = "  -4.78%  "
